$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 9-10; existing rows 9..62 shift down to 11..64.
$ws.Rows("9:10").Insert()

# Fill in the two newly inserted rows with the new data (constant columns match
# the rest of the sheet: Mercado/Mercado ID/Región/Codreg/Tipo/Producto/Categoría).

# Row 9
$ws.Cells.Item(9, 1).Value = 2
$ws.Cells.Item(9, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(9, 3).Value = "Coquimbo"
$ws.Cells.Item(9, 4).Value = 44685
$ws.Cells.Item(9, 5).Value = 4
$ws.Cells.Item(9, 6).Value = "Fruta"
$ws.Cells.Item(9, 7).Value = 100103
$ws.Cells.Item(9, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(9, 9).Value = 100103002
$ws.Cells.Item(9, 10).Value = "Ciruela"
$ws.Cells.Item(9, 11).Value = "Angeleno"
$ws.Cells.Item(9, 12).Value = "Especial"
$ws.Cells.Item(9, 13).Value = 16
$ws.Cells.Item(9, 14).Value = 200000
$ws.Cells.Item(9, 15).Value = 210000
$ws.Cells.Item(9, 16).Value = 205000
$ws.Cells.Item(9, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(9, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(9, 19).Value = 456
$ws.Cells.Item(9, 20).Value = 450

# Row 10
$ws.Cells.Item(10, 1).Value = 2
$ws.Cells.Item(10, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(10, 3).Value = "Coquimbo"
$ws.Cells.Item(10, 4).Value = 44685
$ws.Cells.Item(10, 5).Value = 4
$ws.Cells.Item(10, 6).Value = "Fruta"
$ws.Cells.Item(10, 7).Value = 100103
$ws.Cells.Item(10, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(10, 9).Value = 100103002
$ws.Cells.Item(10, 10).Value = "Ciruela"
$ws.Cells.Item(10, 11).Value = "Angeleno"
$ws.Cells.Item(10, 12).Value = "Primera"
$ws.Cells.Item(10, 13).Value = 20
$ws.Cells.Item(10, 14).Value = 180000
$ws.Cells.Item(10, 15).Value = 190000
$ws.Cells.Item(10, 16).Value = 185000
$ws.Cells.Item(10, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(10, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(10, 19).Value = 411
$ws.Cells.Item(10, 20).Value = 450
